$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Delete column C ("Carbon intensity ANR (kgCO2eq/MWhe)") - all downstream
# columns (D..O) shift left by one (D->C ... O->N), and formulas' cell
# references auto-adjust accordingly.
$ws.Columns.Item(3).Delete()

# The last column (old O, now N) used to hold the formula
# "=C*(D/E)" ("Carbon intensity (kgCO2eq/kgH2)"); after the delete this
# reference is broken (#REF!) because the old C column is gone. Replace it
# with the new literal carbon-intensity figures from the LCA ANL analysis.
$carbonIntensity = @{
    2 = 0.3
    3 = 0.3
    4 = 0.3
    5 = 0.3
    6 = 0.3
    7 = 0.4
    8 = 0.4
    9 = 0.4
    10 = 0.4
    11 = 0.4
    12 = 0.4
    13 = 0.4
    14 = 0.4
    15 = 0.4
    16 = 0.4
}

foreach ($row in $carbonIntensity.Keys) {
    $ws.Range("N$row").Value2 = $carbonIntensity[$row]
}

# Reflect the author's final view state on the sheet (matches the saved
# selection/scroll position recorded in the workbook).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("O13").Select() | Out-Null
